$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.251553416252136
$ws.Range("B1").Value = 2.636661052703857
$ws.Range("C1").Value = 8.29698371887207
$ws.Range("D1").Value = 2.103338003158569
$ws.Range("E1").Value = 1.134089589118958
